$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.989.64'
$ws.Range("E2").Value = '  +3.67%  '

$ws.Range("D3").Value = '1.692.57'
$ws.Range("E3").Value = '  +3.56%  '

$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").Value = "'220.65"
$ws.Range("E5").Value = '  +2.43%  '

$ws.Range("D6").Value = "'0.533"
$ws.Range("E6").Value = '  +2.62%  '

$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("D8").Value = "'29.54"
$ws.Range("E8").Value = '  +2.91%  '

$ws.Range("E9").Value = '  +3.00%  '

$ws.Range("D10").Value = "'0.0639"
$ws.Range("E10").Value = '  +4.86%  '

$ws.Range("D11").Value = "'0.0911"
$ws.Range("E11").Value = '  +1.25%  '

$ws.Range("D12").Value = '1.935.01'
$ws.Range("E12").Value = '  +3.54%  '

$ws.Range("D13").Value = '1.705.86'
$ws.Range("E13").Value = '  +4.39%  '

$ws.Range("D14").Value = "'10.20"
$ws.Range("E14").Value = '  +8.03%  '

$ws.Range("E15").Value = '  +3.44%  '

$ws.Range("E16").Value = '  +6.23%  '

$ws.Range("D17").Value = '31.032.00'
$ws.Range("E17").Value = '  +3.73%  '

$ws.Range("D18").Value = "'66.88"
$ws.Range("E18").Value = '  +3.48%  '

$ws.Range("D19").Value = "'247.57"
$ws.Range("E19").Value = '  +2.97%  '

$ws.Range("D20").Value = '0.0₃0720'
$ws.Range("E20").Value = '  +2.46%  '

$ws.Range("E21").Value = '  +0.35%  '

$ws.Range("D22").Value = "'4.28"
$ws.Range("E22").Value = '  +3.52%  '

$ws.Range("E23").Value = '  +2.03%  '

$ws.Range("E24").Value = '  -0.83%  '

$ws.Range("D25").Value = "'158.47"
$ws.Range("E25").Value = '  +0.68%  '

$ws.Range("E26").Value = '  +2.72%  '

$ws.Range("E27").Value = '  +2.37%  '

$ws.Range("D28").Value = "'6.72"
$ws.Range("E28").Value = '  +1.41%  '

$ws.Range("D30").Value = "'0.0501"
$ws.Range("E30").Value = '  +2.51%  '

$ws.Range("D31").Value = "'3.60"
$ws.Range("E31").Value = '  +6.26%  '

$ws.Range("E32").Value = '  +3.70%  '

$ws.Range("E33").Value = '  +5.10%  '

$ws.Range("D34").Value = '1.519.07'
$ws.Range("E34").Value = '  +6.86%  '

$ws.Range("E35").Value = '  +2.42%  '

$ws.Range("E36").Value = '  +0.85%  '

$ws.Range("E37").Value = '  +10.56%  '

$ws.Range("D38").Value = "'82.79"
$ws.Range("E38").Value = '  +9.01%  '

$ws.Range("E39").Value = '  +4.33%  '

$ws.Range("B40").Value = 'HuobiToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D40").Value = "'2.30"
$ws.Range("E40").Value = '  +0.42%  '

$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = "'2.66"
$ws.Range("E41").Value = '  -4.00%  '

$ws.Range("E42").Value = '  +2.30%  '

$ws.Range("D43").Value = "'0.847"
$ws.Range("E43").Value = '  +1.61%  '

$ws.Range("D44").Value = "'0.0504"
$ws.Range("E44").Value = '  +0.84%  '

$ws.Range("D45").Value = "'1.04"
$ws.Range("E45").Value = '  +2.90%  '

$ws.Range("E46").Value = '  +0.10%  '

$ws.Range("E47").Value = '  +4.47%  '

$ws.Range("D48").Value = "'51.77"
$ws.Range("E48").Value = '  +6.80%  '

$ws.Range("D49").Value = '1.822.77'
$ws.Range("E49").Value = '  +2.61%  '

$ws.Range("E50").Value = '  +8.00%  '

$ws.Range("D51").Value = "'93.93"
$ws.Range("E51").Value = '  +1.01%  '
